$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the table header cells (renumbering / relabeling the columns).
# Updating the header row cell text also updates the underlying table
# column names, since the table header row is bound to A3:I3.
$ws.Range("A3").Value2 = "01_Object"
$ws.Range("B3").Value2 = "02_Type"
$ws.Range("C3").Value2 = "03_Link to Specification"
$ws.Range("D3").Value2 = "04_Developer"
$ws.Range("E3").Value2 = "05_Due Date"
$ws.Range("F3").Value2 = "06_Sprint"
$ws.Range("G3").Value2 = "07_Ready for %"
$ws.Range("H3").Value2 = "08_Dev Comment "
$ws.Range("I3").Value2 = "09_PM Comment"

# Update the selected cell to G3 (matches the recorded selection change)
$ws.Activate()
$ws.Range("G3").Select()
